# Update the "as_of_utc" timestamp column (AA) on the data sheets
# from "2025-11-06 03:04:31" to "2025-11-06 07:04:09" for rows 2-26.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # column AA = 27
        if ($cell.Value() -eq "2025-11-06 03:04:31") {
            $cell.Value = "2025-11-06 07:04:09"
        }
    }
}
